# Update the "Enterprises density (per 1000 people)" and
# "Enterprises (% of total)" rows on the Summary sheet with refreshed
# figures (Micro / SMEs / MSMEs columns).
#
# The values are stored as text (shared strings) in the workbook, not as
# numbers, so we force text entry (leading apostrophe) and then restore
# each cell's original style, since just assigning a numeric-looking
# string would otherwise make Excel coerce the cell into a Number type.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "B11" = "0.37"
    "C11" = "0.04"
    "D11" = "0.41"
    "B12" = "89.01"
    "C12" = "10.13"
    "D12" = "99.14"
}

foreach ($cellRef in @("B11", "C11", "D11", "B12", "C12", "D12")) {
    $range = $ws.Range($cellRef)
    $origStyle = $range.Style
    $range.Value = "'" + $updates[$cellRef]
    $range.Style = $origStyle
}
